$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 570, pushing existing rows 570:601 down to 571:602
$ws.Rows.Item(570).Insert()

# Populate the newly inserted row 570 with the weekly Perejil price record
$ws.Cells.Item(570, 1).Value = 9
$ws.Cells.Item(570, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(570, 3).Value = "Metropolitana"
$ws.Cells.Item(570, 4).Value = 45147
$ws.Cells.Item(570, 5).Value = 13
$ws.Cells.Item(570, 6).Value = 100112044
$ws.Cells.Item(570, 7).Value = "Perejil"
$ws.Cells.Item(570, 8).Value = "Sin especificar"
$ws.Cells.Item(570, 9).Value = "Primera"
$ws.Cells.Item(570, 10).Value = 70
$ws.Cells.Item(570, 11).Value = 15000
$ws.Cells.Item(570, 12).Value = 18000
$ws.Cells.Item(570, 13).Value = 16500
$ws.Cells.Item(570, 14).Value = "`$/docena de atados"
$ws.Cells.Item(570, 15).Value = "Región Metropolitana"
$ws.Cells.Item(570, 16).Value = 5500
$ws.Cells.Item(570, 17).Value = 3
$ws.Cells.Item(570, 18).Value = "Hortaliza"
